$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '258.09'
Set-TextValue 'E2' '3.47%'
Set-TextValue 'D3' '27.33'
Set-TextValue 'E3' '-3.97%'
Set-TextValue 'D4' '5.212'
Set-TextValue 'E4' '-0.82%'
Set-TextValue 'D5' '0.05940'
Set-TextValue 'E5' '3.27%'
Set-TextValue 'D6' '6.705'
Set-TextValue 'E6' '0.67%'
Set-TextValue 'D7' '0.8686'
Set-TextValue 'E7' '0.80%'
Set-TextValue 'D8' '1.013'
Set-TextValue 'E8' '10.25%'
Set-TextValue 'E9' '0.91%'
Set-TextValue 'D10' '0.07196'
Set-TextValue 'E10' '0.42%'
Set-TextValue 'D11' '0.03145'
Set-TextValue 'E11' '0.55%'
Set-TextValue 'D12' '0.09253'
Set-TextValue 'E12' '0.17%'
Set-TextValue 'D13' '0.001550'
Set-TextValue 'E13' '0.98%'
Set-TextValue 'D14' '0.0006069'
Set-TextValue 'E14' '0.82%'
Set-TextValue 'D15' '0.006023'
Set-TextValue 'E15' '0.91%'
Set-TextValue 'D16' '3.491'
Set-TextValue 'E16' '-0.26%'
Set-TextValue 'E17' '1.15%'
Set-TextValue 'E18' '-2.18%'
Set-TextValue 'E19' '0.60%'
Set-TextValue 'D20' '0.03557'
Set-TextValue 'E20' '5.30%'
Set-TextValue 'E21' '-0.49%'
Set-TextValue 'D22' '3.555'
Set-TextValue 'E22' '0.45%'
Set-TextValue 'D23' '0.04277'
Set-TextValue 'E23' '2.48%'
Set-TextValue 'D24' '0.1363'
Set-TextValue 'E24' '-0.99%'
Set-TextValue 'D25' '0.001224'
Set-TextValue 'E25' '0.48%'
Set-TextValue 'D26' '0.004517'
Set-TextValue 'E26' '-10.52%'
Set-TextValue 'E27' '0.01%'
Set-TextValue 'E28' '-22.94%'
Set-TextValue 'D40' '0.03837'
Set-TextValue 'E40' '0.01%'
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue 'D41' '0.006577'
Set-TextValue 'E41' '15.88%'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue 'D42' '0.1104'
Set-TextValue 'E42' '2.09%'
Set-TextValue 'D43' '0.002309'
Set-TextValue 'E43' '5.01%'
Set-TextValue 'D44' '0.01051'
Set-TextValue 'E44' '7.69%'
Set-TextValue 'D45' '0.00005485'
Set-TextValue 'E45' '4.04%'
Set-TextValue 'E46' '0.04%'
Set-TextValue 'E47' '28.55%'
Set-TextValue 'D48' '0.002230'
Set-TextValue 'E48' '2.49%'
Set-TextValue 'E49' '0.04%'
Set-TextValue 'E50' '0.04%'
